$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value recorded one by one
$ws.Range("E8").Value = 63

# New row: 台子 (table)
$ws.Range("G11").Value = "台子"
$ws.Range("H11").Value = 7.5

# New row: 银子 (silver)
$ws.Range("G14").Value = "银子"
$ws.Range("H14").Value = 7.8

# Update the selection to reflect where the user last clicked
$ws.Range("H15").Select()
